$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate the "K" column (column G) values to reflect recomputed
# strike counts (s_vals) instead of the old "Strike#" values.
$kValues = @{
    2  = 6
    3  = 0
    4  = 1
    5  = 2
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 2
    12 = 9
    13 = 2
    14 = 7
    15 = 7
    16 = 1
    17 = 2
    18 = 3
    19 = 4
    20 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
